$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F362").Value = 229208
$ws.Range("F363").Value = 188698
$ws.Range("F364").Value = 168406
$ws.Range("F365").Value = 184982
$ws.Range("F366").Value = 339401
$ws.Range("F367").Value = 766227
$ws.Range("F368").Value = 340311
$ws.Range("G368").Value = 2262
$ws.Range("F369").Value = 234944
$ws.Range("F370").Value = 180720
$ws.Range("F371").Value = 160132
$ws.Range("F373").Value = 350140
$ws.Range("G373").Value = 2379
$ws.Range("F375").Value = 345524
$ws.Range("G375").Value = 1842
$ws.Range("F376").Value = 220850
$ws.Range("F377").Value = 176805
$ws.Range("G377").Value = 1820
$ws.Range("F378").Value = 157102
$ws.Range("G378").Value = 1548
$ws.Range("F379").Value = 179740
$ws.Range("G379").Value = 1618
$ws.Range("F380").Value = 344484
$ws.Range("G380").Value = 2015
$ws.Range("F381").Value = 745207
$ws.Range("G381").Value = 2687
$ws.Range("F383").Value = 220514
$ws.Range("G383").Value = 1762
$ws.Range("F384").Value = 171847
$ws.Range("G384").Value = 1511
$ws.Range("F385").Value = 150720
$ws.Range("F386").Value = 182558
$ws.Range("G386").Value = 1358
$ws.Range("F387").Value = 350981
$ws.Range("G387").Value = 1661
$ws.Range("F388").Value = 729119
$ws.Range("G388").Value = 2197
$ws.Range("F390").Value = 219574
$ws.Range("G390").Value = 1472
$ws.Range("F391").Value = 176796
$ws.Range("G391").Value = 1207
$ws.Range("F392").Value = 220764
$ws.Range("G392").Value = 1212
$ws.Range("F393").Value = 307151
$ws.Range("G393").Value = 1225
$ws.Range("F394").Value = 166017
$ws.Range("G394").Value = 632
$ws.Range("F395").Value = 749743
$ws.Range("G395").Value = 1952
$ws.Range("F396").Value = 164636
$ws.Range("G396").Value = 549
$ws.Range("F397").Value = 108142
$ws.Range("G397").Value = 639
$ws.Range("F398").Value = 298295
$ws.Range("F399").Value = 200387
$ws.Range("G399").Value = 968
$ws.Range("F400").Value = 150131
$ws.Range("F401").Value = 273116
$ws.Range("G401").Value = 933
$ws.Range("F402").Value = 716236
$ws.Range("G402").Value = 1382
$ws.Range("F404").Value = 224347
$ws.Range("G404").Value = 903
$ws.Range("F405").Value = 173594
$ws.Range("F406").Value = 170406
$ws.Range("G406").Value = 676
$ws.Range("F407").Value = 157165
$ws.Range("F408").Value = 301632
$ws.Range("G408").Value = 835
$ws.Range("F409").Value = 692569
$ws.Range("F410").Value = 349168
$ws.Range("G410").Value = 616
$ws.Range("F411").Value = 222845
$ws.Range("G411").Value = 817
$ws.Range("F412").Value = 173835
$ws.Range("G412").Value = 638
$ws.Range("F413").Value = 147350
$ws.Range("G413").Value = 652
